# "version estable 20230214 1540"
#
# Hoja3 ("Hoja3") was an empty placeholder sheet. This edit populates it
# with the same 52-week reading plan that lives on Hoja2, then tweaks
# week 5 (row 6) so Levitico is split into four shorter readings instead
# of the original four-chapter chunks, highlighting the changed cell.
# Hoja3 becomes the active/selected tab (with D3 selected), matching the
# commit's final on-screen state.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Hoja2")
$ws3 = $wb.Worksheets.Item("Hoja3")

# --- Copy the full A1:H53 reading-plan table from Hoja2 onto Hoja3 ---
for ($r = 1; $r -le 53; $r++) {
    for ($c = 1; $c -le 8; $c++) {
        $srcCell = $ws2.Cells.Item($r, $c)
        $dstCell = $ws3.Cells.Item($r, $c)
        $dstCell.Value = $srcCell.Value2
    }
}

# --- Week 5 (row 6): re-split the Levitico readings ---
$ws3.Range("B6").Value = "Levítico (1 a 2)"
$ws3.Range("C6").Value = "Levítico (2 a 4)"
$ws3.Range("D6").Value = "Levítico (5 a 9)"
$ws3.Range("E6").Value = "Levítico (10 a 14)"

# Highlight the edited cell with a yellow fill.
$ws3.Range("D6").Interior.Color = 65535

# --- Selection / active-tab bookkeeping ---
# Hoja2 keeps its old "whole sheet" selection but is no longer the tab
# shown on open; Hoja3 takes over as the selected tab, parked on D3.
$ws2.Range("A1:H1048576").Select()
$ws3.Activate()
$ws3.Range("D3").Select()
